# Add a new BOM line (D36 - TVS diode, 1.5KE16CA) as row 30, shifting the
# existing rows (old row 30 onward) down by one. Excel's native row-insert
# takes care of re-pointing formulas / shared-formula ranges and the
# K2 "SUM(J2:J146)" style range that intentionally overshoots the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 30 (RUEF1600 / F1 entry moves to 31).
$ws.Rows("30:30").Insert()

# Populate the new row with the TVS diode part.
$ws.Cells.Item(30, 1).Value = "D36"              # Part
$ws.Cells.Item(30, 2).Value = "18.8V SO"          # Value
$ws.Cells.Item(30, 3).Value = "TVS"               # Device
$ws.Cells.Item(30, 4).Value = "DO-201"            # Package
$ws.Cells.Item(30, 5).Value = "ME"                # Vendor
$ws.Cells.Item(30, 7).Value = "1.5KE16CA"         # Manufacturer PN (entered before Vendor PN)
$ws.Cells.Item(30, 6).Value = "576-1.5KE16CA"     # Vendor PN

# Match the author's final selection (F31).
$ws.Range("F31").Select()
